# Insert a new weekly data row before the existing row 235.
# This shifts all existing data rows 235-330 down to 236-331,
# growing the sheet's used range from A1:R330 to A1:R331.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("235:235").Insert()

# Populate the newly inserted row 235 with the new record.
$ws.Range("A235").Value = 11
$ws.Range("B235").Value = "Vega Monumental Concepción"
$ws.Range("C235").Value = "Bíobío"
$ws.Range("D235").Value = 45202
$ws.Range("E235").Value = 8
$ws.Range("F235").Value = 100112003
$ws.Range("G235").Value = "Ajo"
$ws.Range("H235").Value = "Chino"
$ws.Range("I235").Value = "Primera"
$ws.Range("J235").Value = 200
$ws.Range("K235").Value = 21000
$ws.Range("L235").Value = 21000
$ws.Range("M235").Value = 21000
$ws.Range("N235").Value = "$/caja 10 kilos"
$ws.Range("O235").Value = "China"
$ws.Range("P235").Value = 2100
$ws.Range("Q235").Value = 10
$ws.Range("R235").Value = "Hortaliza"
